# Final cleanup pass of the GRNmap test-file audit.
# - Removes the stray "Sheet" label row (with placeholder values 3/4) that had
#   been left in the optimization_parameters sheet, shifting the
#   simulation_timepoints row up to take its place.
# - Updates selections left over from interactive review.
# - Leaves threshold_b as the active/selected sheet.

$wb = $excel.ActiveWorkbook

$wsProd  = $wb.Worksheets.Item("production_rates")
$wsOpt   = $wb.Worksheets.Item("optimization_parameters")
$wsThr   = $wb.Worksheets.Item("threshold_b")

# Delete the obsolete "Sheet" row (row 16) from optimization_parameters.
$wsOpt.Rows.Item(16).Delete()

# Restore the selection on production_rates left over from scrolling.
$wsProd.Range("C40").Select()

# Select the row that now occupies row 16 (formerly row 17,
# simulation_timepoints) on optimization_parameters.
$wsOpt.Rows.Item(16).Select()

# Make threshold_b the active sheet/tab, with its prior cell selection.
$wsThr.Activate()
$wsThr.Range("A2").Select()
